# Apply the "Append: 2025-10-03 06:26 JST" update to the ランサーズ sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newDate = "2025-10-03 06:26:20"

# --- Remove old hyperlinks (the engine only supports clearing the whole
#     collection at once; we re-add the ones we keep further down). ---
$ws.Hyperlinks.Delete()

# --- Row 2 : timestamp refresh only ---
$ws.Range("A2").Value = $newDate

# --- Row 3 : new listing ---
$ws.Range("A3").Value = $newDate
$ws.Range("B3").Value = "【業務自動化×補助金対応】生成AI活用/日本人モデル画像生成歓迎"
$ws.Range("D3").Value = "3,000,000 円 ~ 5,000,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5405834"
$ws.Range("G3").Value = 395
$ws.Range("H3").Value = "🔥AI,Ai ◆自動化"

# --- Row 4 : new listing ---
$ws.Range("A4").Value = $newDate
$ws.Range("B4").Value = "Reactの細かい修正の対応"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5405740"
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = "🔥React"

# --- Row 5 : new listing ---
$ws.Range("A5").Value = $newDate
$ws.Range("B5").Value = "【急募】愛知県でのBtoB受発注システム開発者募集"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5405971"
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = "◆開発,システム開発"

# --- Row 6 : new listing ---
$ws.Range("A6").Value = $newDate
$ws.Range("B6").Value = "【WEB】Nuxt3でのWEBページ表示速度改善、他継続して弊社システムの開発に携われる方"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5406001"
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = "◆開発"

# --- Row 7 : new listing ---
$ws.Range("A7").Value = $newDate
$ws.Range("B7").Value = "【急募】Excelマクロでデータからグラフを自動作成するツール"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5405961"
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = "◆ツール"

# --- Row 8 : new listing, skill-summary column now empty ---
$ws.Range("A8").Value = $newDate
$ws.Range("B8").Value = "【急募】国内300店舗規模のスーパーマーケット向けActive Directory構築"
$ws.Range("D8").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5406008"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# --- Row 9 : new listing, skill-summary column now empty ---
$ws.Range("A9").Value = $newDate
$ws.Range("B9").Value = "【急募】全国物件情報抽出プログラム作成依頼"
$ws.Range("D9").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5405763"
$ws.Range("G9").Value = 13
$ws.Range("H9").ClearContents()

# --- Drop the old rows 10-17 entirely ---
$ws.Rows("10:17").Delete()

# --- Re-create the hyperlinks for the rows that remain (F2:F9) ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405813")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5405834")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5405740")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5405971")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5406001")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5405961")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5406008")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5405763")

# --- Column width tweaks (COM ColumnWidth units run 5/6 low vs. the
#     stored OOXML width, so compensate before assigning). ---
$ws.Columns("B").ColumnWidth = 48 - 5/6
$ws.Columns("D").ColumnWidth = 32 - 5/6
$ws.Columns("H").ColumnWidth = 19 - 5/6
